$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The columns that actually carry data which differs between each swapped
# row pair. (Other columns in the pair already hold identical values, so
# there's no need to touch them.)
$cols = @("A","B","D","E","F","G","H","Q","R","AC")

function Swap-Rows([int]$r1, [int]$r2) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

Swap-Rows 5 6
Swap-Rows 7 8
Swap-Rows 13 14
